$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in weather data for the existing cities (rows 2-10)
$ws.Range("B2").Value = 46.65
$ws.Range("C2").Value = 39.79
$ws.Range("D2").Value = 34

$ws.Range("B3").Value = 58.73
$ws.Range("C3").Value = 58.73
$ws.Range("D3").Value = 94

$ws.Range("B4").Value = 51.71
$ws.Range("C4").Value = 49.96
$ws.Range("D4").Value = 72

$ws.Range("B5").Value = 10.87
$ws.Range("C5").Value = -1.73
$ws.Range("D5").Value = 67

$ws.Range("B6").Value = 66.22
$ws.Range("C6").Value = 66.11
$ws.Range("D6").Value = 76

$ws.Range("B7").Value = 54.72
$ws.Range("C7").Value = 54.5
$ws.Range("D7").Value = 98

$ws.Range("B8").Value = 50.29
$ws.Range("C8").Value = 49.15
$ws.Range("D8").Value = 88

$ws.Range("B9").Value = 35.44
$ws.Range("C9").Value = 28.99
$ws.Range("D9").Value = 60

$ws.Range("B10").Value = 51.48
$ws.Range("C10").Value = 49.28
$ws.Range("D10").Value = 63

# New cities' names added first (Tampa Bay, then Midland) to keep shared-string order stable
$ws.Range("A11").Value = "Tampa Bay"
$ws.Range("A12").Value = "Midland"

# New city "Midland" - successful lookup
$ws.Range("B12").Value = 45.73
$ws.Range("C12").Value = 40.229999999999997
$ws.Range("D12").Value = 48

# New city "Tampa Bay" - simulate a lookup failure (control flow / error handling)
$ws.Range("B11").Value = "Not Found"
$ws.Range("C11").Value = "Not Found"
$ws.Range("D11").Value = "Not Found"

# Expand the table to include the two new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D12"))

# Update the active selection to match the author's final cursor position
$ws.Range("J7").Select()
